$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("matrizsod")
$ws.Range("A1").Value = "X"
$ws.Activate()
$ws.Range("A8").Select()
